# Weekly update: insert a new pair of rows (Camote 1a/2a nueva(o)) at the
# top of the "Zapallo / Camote" block (rows 316-317), pushing the existing
# data down by two rows. The sheet's used range grows from A1:R345 to
# A1:R347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 316 downwards (and below) by two rows to make room for the
# new weekly records.
$ws.Rows("316:317").Insert()

# --- New row 316 --------------------------------------------------------
$ws.Cells.Item(316, 1).Value = 8
$ws.Cells.Item(316, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(316, 3).Value = "Coquimbo"
$ws.Cells.Item(316, 4).Value = 44449
$ws.Cells.Item(316, 5).Value = 4
$ws.Cells.Item(316, 6).Value = 100112045
$ws.Cells.Item(316, 7).Value = "Zapallo"
$ws.Cells.Item(316, 8).Value = "Camote"
$ws.Cells.Item(316, 9).Value = "1a nueva(o)"
$ws.Cells.Item(316, 10).Value = 800
$ws.Cells.Item(316, 11).Value = 950
$ws.Cells.Item(316, 12).Value = 1000
$ws.Cells.Item(316, 13).Value = 975
$ws.Cells.Item(316, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(316, 15).Value = "Perú"
$ws.Cells.Item(316, 16).Value = 975
$ws.Cells.Item(316, 17).Value = 1
$ws.Cells.Item(316, 18).Value = "Hortaliza"

# --- New row 317 --------------------------------------------------------
$ws.Cells.Item(317, 1).Value = 8
$ws.Cells.Item(317, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(317, 3).Value = "Coquimbo"
$ws.Cells.Item(317, 4).Value = 44449
$ws.Cells.Item(317, 5).Value = 4
$ws.Cells.Item(317, 6).Value = 100112045
$ws.Cells.Item(317, 7).Value = "Zapallo"
$ws.Cells.Item(317, 8).Value = "Camote"
$ws.Cells.Item(317, 9).Value = "2a nueva(o)"
$ws.Cells.Item(317, 10).Value = 560
$ws.Cells.Item(317, 11).Value = 850
$ws.Cells.Item(317, 12).Value = 900
$ws.Cells.Item(317, 13).Value = 875
$ws.Cells.Item(317, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(317, 15).Value = "Perú"
$ws.Cells.Item(317, 16).Value = 875
$ws.Cells.Item(317, 17).Value = 1
$ws.Cells.Item(317, 18).Value = "Hortaliza"
